# Swap the order of the first two comma-separated entries in the
# "Recorded By" column (G) for the known recorder-name/email combinations.
# This mirrors a bulk find/replace across the sheet where 4 distinct
# whole-cell text values were replaced with their "swapped" counterparts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$replacements = @{
    "System, backup@backdoor.com, system" = "backup@backdoor.com, System, system"
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com"
    "System, backup@backdoor.com"         = "backup@backdoor.com, System"
    "admin@admin.com, System"             = "System, admin@admin.com"
}

$lastRow = $ws.UsedRange.Rows.Count
$changed = 0

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $txt = $cell.Text
    if ($replacements.ContainsKey($txt)) {
        $cell.Value = $replacements[$txt]
        $changed = $changed + 1
    }
}

Write-Host ("Rows updated: " + $changed)
